$wb = $excel.ActiveWorkbook

# --- Sheet "Balraj": new daily-tracker rows for 2020-10-05 / 2020-10-06 ---
$ws1 = $wb.Worksheets.Item("Balraj")

# Row 4 (entry #3, RPA SAMSUNG)
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = (Get-Date -Year 2020 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws1.Range("B4").NumberFormat = "mm-dd-yy"
$ws1.Range("C4").Value = "RPA SAMSUNG"
$ws1.Range("D4").Value = "1. DRS-Weekly correction received and completed"
$ws1.Range("E4").Value = 1
$ws1.Range("E4").NumberFormat = "0%"
$ws1.Range("F4").Value = "Completed"

# Row 5 (continuation of entry #3)
$ws1.Range("D5").Value = "2. Return Credit Correction received and completed "
$ws1.Range("E5").Value = 1
$ws1.Range("E5").NumberFormat = "0%"
$ws1.Range("F5").Value = "Completed"

# Row 6 (entry #4, RPA SONY)
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = (Get-Date -Year 2020 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws1.Range("B6").NumberFormat = "mm-dd-yy"
$ws1.Range("C6").Value = "RPA SONY"
$ws1.Range("D6").Value = "Daily Task of Scheduling Report has been completed for download and upload (11 files `nexcept one is having upload issue which under progress by Mohan san)"
$ws1.Range("D6").WrapText = $true
$ws1.Range("E6").Value = 1
$ws1.Range("E6").NumberFormat = "0%"
$ws1.Range("F6").Value = "Completed"

# Row 7 stays blank (untouched)

# Row 8 (entry #5, RPA SAMSUNG)
$ws1.Range("A8").Value = 5
$ws1.Range("B8").Value = (Get-Date -Year 2020 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws1.Range("B8").NumberFormat = "mm-dd-yy"
$ws1.Range("C8").Value = "RPA SAMSUNG"
$ws1.Range("D8").Value = "1. DRS-WeeKly task of SSC8 , Correction has been done"
$ws1.Range("E8").Value = 1
$ws1.Range("E8").NumberFormat = "0%"
$ws1.Range("F8").Value = "Completed"

# Row 9 (continuation of entry #5)
$ws1.Range("D9").Value = "2. DRS-DAILY task of SSC11, Correction has been done"
$ws1.Range("E9").Value = 1
$ws1.Range("E9").NumberFormat = "0%"
$ws1.Range("F9").Value = "Completed"

# Row 10 (entry #6, RPA SONY)
$ws1.Range("A10").Value = 6
$ws1.Range("B10").Value = (Get-Date -Year 2020 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws1.Range("B10").NumberFormat = "mm-dd-yy"
$ws1.Range("C10").Value = "RPA SONY"
$ws1.Range("D10").Value = "1. SchedulingReport daiky has been completed"
$ws1.Range("E10").Value = 1
$ws1.Range("E10").NumberFormat = "0%"
$ws1.Range("F10").Value = "Completed"

# Make "Balraj" the active/selected sheet and restore its selection to A2:F10.
# Activating it automatically clears tabSelected on the previously-active
# sheet ("Monisha"), matching the diff without touching it directly.
$ws1.Activate()
$ws1.Range("A2:F10").Select()
